# Updating escape data preprocessing to fix lab units
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix lab unit cells (g/L -> g/dL typos, fill in missing unit columns) ---

# ALB (row 95): ESCAPE unit column B was mislabeled g/L, should be g/dL
$ws.Range("B95").Value = "g/dL"

# BUN (row 98) / CRT (row 99): add the missing Cardiogenic Shock (F) and
# Serial Cardiac (G) unit columns, matching the existing mg/dL columns
$ws.Range("F98").Value = "mg/dL"
$ws.Range("G98").Value = "mg/dL"
$ws.Range("F99").Value = "mg/dL"
$ws.Range("G99").Value = "mg/dL"

# HEM (row 102): ESCAPE unit column B was mislabeled g/L, should be g/dL
$ws.Range("B102").Value = "g/dL"

# TOTP (row 107): ESCAPE unit column B was mislabeled g/L, should be g/dL
$ws.Range("B107").Value = "g/dL"

# SOD (row 105): add the missing Cardiogenic Shock (F) and Serial Cardiac
# (G) unit columns, matching the existing mmol/L columns
$ws.Range("F105").Value = "mmol/L"
$ws.Range("G105").Value = "mmol/L"

# WBC (row 108): add the missing Cardiogenic Shock (F) and Serial Cardiac
# (G) unit columns
$ws.Range("F108").Value = "10^9/L"
$ws.Range("G108").Value = "10^9/L"

# PLA (row 103) / POT (row 104): note that the ESCAPE/GUIDE-IT units are the
# same as BEST's, flagged in a new column E
$ws.Range("E103").Value = "<- note these same"
$ws.Range("E104").Value = "<- note these same"

# --- Freeze header rows/column on the sheet view ---
$ws.Range("B5").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B18").Select() | Out-Null
